# Added timestamp to code
# - Inserts a new "date" column in column A (header "date" in A1), with a
#   YYYY-MM-DD HH:MM:SS timestamp for every data row.
# - Refreshes several of the account-holdings data values (balances,
#   available, holds, prices, dollar values) for the latest pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = 44382.89266203704

# --- New "date" header cell (A1), formatted like the rest of the header row ---
$ws.Range("A1").Value = "date"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats - copy B1's header style onto A1

# --- Column A: replace the old row-index values with the timestamp ---
$ws.Range("A2:A8").Value2 = $timestamp

# Apply the header's base style to the timestamp cells first (border/font/
# alignment), then layer in the custom date/time number format.
$ws.Range("B1").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Refresh data values (balances / available / holds / price / dollar value) ---

# Row 3 - USDT
$ws.Range("D3").Value = 1011.67532122
$ws.Range("E3").Value = 120.41278248
$ws.Range("F3").Value = 891.26253874
$ws.Range("H3").Value = 1011.68

# Row 4 - BTC
$ws.Range("G4").Value = 34055
$ws.Range("H4").Value = 524.5

# Row 5 - ATOM
$ws.Range("D5").Value = 13.1249
$ws.Range("F5").Value = 12.8877
$ws.Range("G5").Value = 13.8091
$ws.Range("H5").Value = 181.24

# Row 6 - ALGO
$ws.Range("G6").Value = 0.8855

# Row 7 - ETH
$ws.Range("G7").Value = 2231.28
